# Auto-generated edit script: updates Price (D) and Volume(1h) (E) columns
# for the cryptos worksheet, preserving text formatting of numeric-looking values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "36.408.37"
$ws.Range("E2").Value = "  -2.52%  "

Set-TextValue $ws.Range("D3") "1.986.59"
$ws.Range("E3").Value = "  -1.31%  "

$ws.Range("E4").Value = "  +0.00%  "

Set-TextValue $ws.Range("D5") "235.14"
$ws.Range("E5").Value = "  -8.93%  "

Set-TextValue $ws.Range("D6") "0.595"
$ws.Range("E6").Value = "  -3.42%  "

Set-TextValue $ws.Range("D8") "54.41"
$ws.Range("E8").Value = "  -3.34%  "

$ws.Range("E9").Value = "  -5.30%  "

Set-TextValue $ws.Range("D10") "57.57"
$ws.Range("E10").Value = "  +2.05%  "

$ws.Range("E11").Value = "  -3.19%  "

Set-TextValue $ws.Range("D12") "0.0987"
$ws.Range("E12").Value = "  -2.96%  "

Set-TextValue $ws.Range("D13") "14.12"

Set-TextValue $ws.Range("D14") "2.277.58"
$ws.Range("E14").Value = "  -1.38%  "

Set-TextValue $ws.Range("D15") "20.12"
$ws.Range("E15").Value = "  -3.44%  "

Set-TextValue $ws.Range("D16") "0.753"
$ws.Range("E16").Value = "  -6.55%  "

Set-TextValue $ws.Range("D17") "5.05"
$ws.Range("E17").Value = "  -3.79%  "

Set-TextValue $ws.Range("D18") "1.991.65"
$ws.Range("E18").Value = "  -1.45%  "

Set-TextValue $ws.Range("D19") "36.370.37"
$ws.Range("E19").Value = "  -2.39%  "

Set-TextValue $ws.Range("D20") "67.61"
$ws.Range("E20").Value = "  -3.10%  "

Set-TextValue $ws.Range("D21") "0.0₃0803"
$ws.Range("E21").Value = "  -4.19%  "

Set-TextValue $ws.Range("D23") "221.12"
$ws.Range("E23").Value = "  -3.15%  "

$ws.Range("E24").Value = "  -0.05%  "

$ws.Range("E25").Value = "  -0.07%  "

$ws.Range("E26").Value = "  -9.44%  "

Set-TextValue $ws.Range("D27") "162.53"
$ws.Range("E27").Value = "  -1.41%  "

Set-TextValue $ws.Range("D28") "8.66"
$ws.Range("E28").Value = "  -3.99%  "

Set-TextValue $ws.Range("D29") "0.129"
$ws.Range("E29").Value = "  -0.72%  "

$ws.Range("E30").Value = "  -4.67%  "

$ws.Range("E31").Value = "  -0.92%  "

Set-TextValue $ws.Range("D32") "0.116"
$ws.Range("E32").Value = "  -3.21%  "

Set-TextValue $ws.Range("D33") "4.35"
$ws.Range("E33").Value = "  -6.50%  "

$ws.Range("E34").Value = "  -7.03%  "

$ws.Range("E35").Value = "  -7.33%  "

$ws.Range("E36").Value = "  -2.76%  "

$ws.Range("E37").Value = "  +0.03%  "

Set-TextValue $ws.Range("D38") "3.34"
$ws.Range("E38").Value = "  -0.94%  "

Set-TextValue $ws.Range("D39") "1.76"
$ws.Range("E39").Value = "  -3.28%  "

Set-TextValue $ws.Range("D40") "5.52"
$ws.Range("E40").Value = "  +5.46%  "

$ws.Range("E41").Value = "  -0.83%  "

Set-TextValue $ws.Range("D42") "1.453.23"
$ws.Range("E42").Value = "  +4.24%  "

Set-TextValue $ws.Range("D43") "0.0931"
$ws.Range("E43").Value = "  +0.08%  "

Set-TextValue $ws.Range("D44") "0.0202"
$ws.Range("E44").Value = "  -5.25%  "

Set-TextValue $ws.Range("D45") "1.09"
$ws.Range("E45").Value = "  -9.37%  "

Set-TextValue $ws.Range("D46") "89.07"
$ws.Range("E46").Value = "  -1.12%  "

Set-TextValue $ws.Range("D47") "14.99"
$ws.Range("E47").Value = "  -4.84%  "

$ws.Range("E48").Value = "  -3.59%  "

$ws.Range("E49").Value = "  -0.88%  "

Set-TextValue $ws.Range("D50") "6.83"
$ws.Range("E50").Value = "  -4.05%  "

Set-TextValue $ws.Range("D51") "3.66"
$ws.Range("E51").Value = "  +6.11%  "
